$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 1464.3518
$ws_ALC.Range("J17").Value = 1515.2439
$ws_ALC.Range("L17").Value = 4545.7317
$ws_ALC.Range("N17").Value = -4881.7317

# ALC row 51
$ws_ALC.Range("H51").Value = 4155.357
$ws_ALC.Range("J51").Value = 4531.1665
$ws_ALC.Range("L51").Value = 4531.1665
$ws_ALC.Range("N51").Value = -5499.1665

# ALC row 54
$ws_ALC.Range("H54").Value = 39683.332
$ws_ALC.Range("J54").Value = 39683.332
$ws_ALC.Range("L54").Value = 39683.332
$ws_ALC.Range("N54").Value = -40655.332

# ALC row 126
$ws_ALC.Range("H126").Value = 41977
$ws_ALC.Range("J126").Value = 41977
$ws_ALC.Range("L126").Value = 41977
$ws_ALC.Range("N126").Value = -51857

# ALC row 132
$ws_ALC.Range("H132").Value = 105212.8
$ws_ALC.Range("I132").Value = 116122.83
$ws_ALC.Range("J132").Value = 9750
$ws_ALC.Range("K132").Value = 348368.49
$ws_ALC.Range("L132").Value = 29250
$ws_ALC.Range("M132").Value = -345838.49
$ws_ALC.Range("N132").Value = -34310

# ALC row 141
$ws_ALC.Range("H141").Value = 2077.75
$ws_ALC.Range("I141").Value = 1891.08
$ws_ALC.Range("K141").Value = 5673.24
$ws_ALC.Range("M141").Value = -493.2399999999998

# ARM row 2
$ws_ARM.Range("H2").Value = 598.06665
$ws_ARM.Range("I2").Value = 520.5789
$ws_ARM.Range("J2").Value = 731.9091
$ws_ARM.Range("K2").Value = 520.5789
$ws_ARM.Range("L2").Value = 731.9091
$ws_ARM.Range("M2").Value = -407.5789
$ws_ARM.Range("N2").Value = -957.9091

# ARM row 110
$ws_ARM.Range("H110").Value = 699.74286
$ws_ARM.Range("I110").Value = 663.3929000000001
$ws_ARM.Range("J110").Value = 845.1429000000001
$ws_ARM.Range("K110").Value = 663.3929000000001
$ws_ARM.Range("L110").Value = 845.1429000000001
$ws_ARM.Range("M110").Value = 1381.6071
$ws_ARM.Range("N110").Value = -4935.1429

# ARM row 116
$ws_ARM.Range("H116").Value = 598.06665
$ws_ARM.Range("I116").Value = 520.5789
$ws_ARM.Range("J116").Value = 731.9091
$ws_ARM.Range("K116").Value = 520.5789
$ws_ARM.Range("L116").Value = 731.9091
$ws_ARM.Range("M116").Value = 1773.4211
$ws_ARM.Range("N116").Value = -5319.9091

# ARM row 122
$ws_ARM.Range("H122").Value = 1939.0968
$ws_ARM.Range("I122").Value = 1090.6316
$ws_ARM.Range("J122").Value = 3282.5
$ws_ARM.Range("K122").Value = 3271.8948
$ws_ARM.Range("L122").Value = 9847.5
$ws_ARM.Range("M122").Value = -821.8948
$ws_ARM.Range("N122").Value = -14747.5

# ARM row 137
$ws_ARM.Range("H137").Value = 42939.6
$ws_ARM.Range("J137").Value = 42939.6
$ws_ARM.Range("L137").Value = 42939.6
$ws_ARM.Range("N137").Value = -53139.6

# BSM row 3
$ws_BSM.Range("H3").Value = 598.06665
$ws_BSM.Range("I3").Value = 520.5789
$ws_BSM.Range("J3").Value = 731.9091
$ws_BSM.Range("K3").Value = 520.5789
$ws_BSM.Range("L3").Value = 731.9091
$ws_BSM.Range("M3").Value = -406.5789
$ws_BSM.Range("N3").Value = -959.9091

# BSM row 20
$ws_BSM.Range("H20").Value = 1761.5
$ws_BSM.Range("I20").Value = 1223.6666
$ws_BSM.Range("K20").Value = 1223.6666
$ws_BSM.Range("M20").Value = -976.6666

# BSM row 107
$ws_BSM.Range("H107").Value = 1068.4166
$ws_BSM.Range("I107").Value = 983.7273
$ws_BSM.Range("K107").Value = 983.7273
$ws_BSM.Range("M107").Value = 936.2727

# BSM row 122
$ws_BSM.Range("H122").Value = 41783.332
$ws_BSM.Range("J122").Value = 41783.332
$ws_BSM.Range("L122").Value = 41783.332
$ws_BSM.Range("N122").Value = -51583.332

# BSM row 137
$ws_BSM.Range("H137").Value = 41308
$ws_BSM.Range("J137").Value = 41308
$ws_BSM.Range("L137").Value = 41308
$ws_BSM.Range("N137").Value = -51508

# CRP row 31
$ws_CRP.Range("H31").Value = 50006510
$ws_CRP.Range("I31").Value = 0
$ws_CRP.Range("J31").Value = 50006510
$ws_CRP.Range("K31").Value = 0
$ws_CRP.Range("L31").Value = 50006510
$ws_CRP.Range("N31").Value = -50007100
$ws_CRP.Range("M31").ClearContents()

# CRP row 34
$ws_CRP.Range("H34").Value = 50006510
$ws_CRP.Range("I34").Value = 0
$ws_CRP.Range("J34").Value = 50006510
$ws_CRP.Range("K34").Value = 0
$ws_CRP.Range("L34").Value = 50006510
$ws_CRP.Range("N34").Value = -50006914
$ws_CRP.Range("M34").ClearContents()

# CRP row 99
$ws_CRP.Range("H99").Value = 14290610
$ws_CRP.Range("I99").Value = 33336404
$ws_CRP.Range("J99").Value = 6264.25
$ws_CRP.Range("K99").Value = 33336404
$ws_CRP.Range("L99").Value = 6264.25
$ws_CRP.Range("M99").Value = -33334906
$ws_CRP.Range("N99").Value = -9260.25

# CRP row 126
$ws_CRP.Range("H126").Value = 14290610
$ws_CRP.Range("I126").Value = 33336404
$ws_CRP.Range("J126").Value = 6264.25
$ws_CRP.Range("K126").Value = 100009212
$ws_CRP.Range("L126").Value = 18792.75
$ws_CRP.Range("M126").Value = -100006742
$ws_CRP.Range("N126").Value = -23732.75

# CRP row 134
$ws_CRP.Range("H134").Value = 8695.25
$ws_CRP.Range("I134").Value = 8778.77
$ws_CRP.Range("J134").Value = 8333.333000000001
$ws_CRP.Range("K134").Value = 26336.31
$ws_CRP.Range("L134").Value = 24999.999
$ws_CRP.Range("M134").Value = -23801.31
$ws_CRP.Range("N134").Value = -30069.999

# CUL row 14
$ws_CUL.Range("H14").Value = 527.52
$ws_CUL.Range("I14").Value = 527.52
$ws_CUL.Range("K14").Value = 1582.56
$ws_CUL.Range("M14").Value = -1409.56

# CUL row 75
$ws_CUL.Range("H75").Value = 3043
$ws_CUL.Range("I75").Value = 313
$ws_CUL.Range("J75").Value = 3589
$ws_CUL.Range("K75").Value = 939
$ws_CUL.Range("L75").Value = 10767
$ws_CUL.Range("M75").Value = 59
$ws_CUL.Range("N75").Value = -12763

# CUL row 78
$ws_CUL.Range("H78").Value = 3043
$ws_CUL.Range("I78").Value = 313
$ws_CUL.Range("J78").Value = 3589
$ws_CUL.Range("K78").Value = 2817
$ws_CUL.Range("L78").Value = 32301
$ws_CUL.Range("M78").Value = 2175
$ws_CUL.Range("N78").Value = -42285

# CUL row 86
$ws_CUL.Range("H86").Value = 2167.5
$ws_CUL.Range("J86").Value = 3126.25
$ws_CUL.Range("L86").Value = 9378.75
$ws_CUL.Range("N86").Value = -11750.75

# CUL row 89
$ws_CUL.Range("H89").Value = 2167.5
$ws_CUL.Range("J89").Value = 3126.25
$ws_CUL.Range("L89").Value = 28136.25
$ws_CUL.Range("N89").Value = -39992.25

# CUL row 116
$ws_CUL.Range("H116").Value = 4998
$ws_CUL.Range("J116").Value = 4998
$ws_CUL.Range("L116").Value = 14994
$ws_CUL.Range("N116").Value = -21878

# CUL row 121
$ws_CUL.Range("H121").Value = 2568.776
$ws_CUL.Range("J121").Value = 2835.577
$ws_CUL.Range("L121").Value = 8506.731
$ws_CUL.Range("N121").Value = -11126.731

# CUL row 131
$ws_CUL.Range("H131").Value = 1035.0754
$ws_CUL.Range("I131").Value = 3616.6667
$ws_CUL.Range("J131").Value = 880.1799999999999
$ws_CUL.Range("K131").Value = 10850.0001
$ws_CUL.Range("L131").Value = 2640.54
$ws_CUL.Range("M131").Value = -5810.000100000001
$ws_CUL.Range("N131").Value = -12720.54

# CUL row 132
$ws_CUL.Range("H132").Value = 2959.0557
$ws_CUL.Range("I132").Value = 959.9
$ws_CUL.Range("J132").Value = 5458
$ws_CUL.Range("K132").Value = 8639.1
$ws_CUL.Range("L132").Value = 49122
$ws_CUL.Range("M132").Value = -6109.1
$ws_CUL.Range("N132").Value = -54182

# CUL row 134
$ws_CUL.Range("H134").Value = 3539.3823
$ws_CUL.Range("I134").Value = 2706.6191
$ws_CUL.Range("J134").Value = 4884.615
$ws_CUL.Range("K134").Value = 8119.8573
$ws_CUL.Range("L134").Value = 14653.845
$ws_CUL.Range("M134").Value = -3049.8573
$ws_CUL.Range("N134").Value = -24793.845

# CUL row 137
$ws_CUL.Range("H137").Value = 3328.05
$ws_CUL.Range("I137").Value = 1236.4445
$ws_CUL.Range("J137").Value = 5039.364
$ws_CUL.Range("K137").Value = 3709.3335
$ws_CUL.Range("L137").Value = 15118.092
$ws_CUL.Range("M137").Value = 1390.6665
$ws_CUL.Range("N137").Value = -25318.092

# GSM row 102
$ws_GSM.Range("H102").Value = 2242.8484
$ws_GSM.Range("I102").Value = 1780.875
$ws_GSM.Range("J102").Value = 3474.7778
$ws_GSM.Range("K102").Value = 1780.875
$ws_GSM.Range("L102").Value = 3474.7778
$ws_GSM.Range("M102").Value = -158.875
$ws_GSM.Range("N102").Value = -6718.7778

# GSM row 137
$ws_GSM.Range("H137").Value = 74324.25
$ws_GSM.Range("J137").Value = 74324.25
$ws_GSM.Range("L137").Value = 74324.25
$ws_GSM.Range("N137").Value = -84524.25

# WVR row 41
$ws_WVR.Range("H41").Value = 8800
$ws_WVR.Range("J41").Value = 8800
$ws_WVR.Range("L41").Value = 8800
$ws_WVR.Range("N41").Value = -9580

# WVR row 126
$ws_WVR.Range("H126").Value = 3253.6
$ws_WVR.Range("I126").Value = 2484.7144
$ws_WVR.Range("J126").Value = 4232.1816
$ws_WVR.Range("K126").Value = 7454.1432
$ws_WVR.Range("L126").Value = 12696.5448
$ws_WVR.Range("M126").Value = -4984.1432
$ws_WVR.Range("N126").Value = -17636.5448

# WVR row 132
$ws_WVR.Range("H132").Value = 13335608
$ws_WVR.Range("I132").Value = 857.9231
$ws_WVR.Range("K132").Value = 2573.7693
$ws_WVR.Range("M132").Value = -43.76929999999993

